# Apply updated cryptocurrency price/volume figures (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.654.52"
$ws.Range("E2").Value = "  +1.66%  "
$ws.Range("D3").Value = "2.491.10"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'492.03"
$ws.Range("E5").Value = "  +2.56%  "
$ws.Range("D6").Value = "'151.34"
$ws.Range("E6").Value = "  +8.89%  "
$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("E8").Value = "  +0.61%  "
$ws.Range("D9").Value = "2.502.06"
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("E10").Value = "  +6.07%  "
$ws.Range("D11").Value = "'0.0983"
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("D12").Value = "'0.337"
$ws.Range("E12").Value = "  +3.14%  "
$ws.Range("D14").Value = "2.922.07"
$ws.Range("E14").Value = "  -0.21%  "
$ws.Range("D15").Value = "56.701.93"
$ws.Range("E15").Value = "  +1.67%  "
$ws.Range("D16").Value = "'21.22"
$ws.Range("E16").Value = "  +3.66%  "
$ws.Range("D17").Value = "'0.0000137"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").Value = "2.499.44"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").Value = "'4.54"
$ws.Range("E19").Value = "  +4.34%  "
$ws.Range("D20").Value = "'10.27"
$ws.Range("E20").Value = "  +3.13%  "
$ws.Range("D21").Value = "'322.66"
$ws.Range("E21").Value = "  +0.71%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("D23").Value = "'5.91"
$ws.Range("E23").Value = "  +4.51%  "
$ws.Range("D24").Value = "'58.80"
$ws.Range("E24").Value = "  +1.60%  "
$ws.Range("E25").Value = "  +1.70%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("D28").Value = "2.598.88"
$ws.Range("E28").Value = "  -0.31%  "
$ws.Range("E29").Value = "  +3.75%  "
$ws.Range("D30").Value = "0.0₃0801"
$ws.Range("E30").Value = "  +4.12%  "
$ws.Range("D31").Value = "'0.998"
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("D32").Value = "'151.00"
$ws.Range("E32").Value = "  +1.33%  "
$ws.Range("D33").Value = "'18.39"
$ws.Range("E33").Value = "  +1.57%  "
$ws.Range("E34").Value = "  +3.21%  "
$ws.Range("D35").Value = "'5.22"
$ws.Range("E35").Value = "  +0.62%  "
$ws.Range("E36").Value = "  +4.81%  "
$ws.Range("E37").Value = "  +2.58%  "
$ws.Range("D38").Value = "'0.875"
$ws.Range("E38").Value = "  +3.52%  "
$ws.Range("E39").Value = "  +5.81%  "
$ws.Range("D40").Value = "'33.97"
$ws.Range("E40").Value = "  -0.85%  "
$ws.Range("D41").Value = "'3.51"
$ws.Range("E41").Value = "  +4.01%  "
$ws.Range("D42").Value = "'0.0558"
$ws.Range("E42").Value = "  +1.66%  "
$ws.Range("D43").Value = "'0.613"
$ws.Range("E43").Value = "  +0.27%  "
$ws.Range("E44").Value = "  -0.42%  "
$ws.Range("E45").Value = "  +8.52%  "
$ws.Range("D46").Value = "'264.05"
$ws.Range("E46").Value = "  +5.53%  "
$ws.Range("E47").Value = "  +2.21%  "
$ws.Range("E48").Value = "  +3.69%  "
$ws.Range("D49").Value = "'10.20"
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("D50").Value = "'17.82"
$ws.Range("E50").Value = "  +2.22%  "
$ws.Range("D51").Value = "1.905.68"
$ws.Range("E51").Value = "  -3.14%  "
